# Update countries & provincias Spain
#
# This script applies the data refresh that:
#   1) Re-orders two countries in the table (Botsuana now sorts right after
#      Burkina Faso / before Niger; Montserrat now sorts right after
#      San Bartolome / before Islas Malvinas), and
#   2) Refreshes the numeric statistics for a number of countries/rows, and
#   3) Bumps the "last updated" timestamp shown in cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Timestamp banner in A1
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 23:59"

# ---------------------------------------------------------------------------
# 2) Plain numeric refreshes (country stays on the same row, only the stats
#    change): row => @(TotalCasos, NuevosCasos, CasosActivos, Recuperados,
#                      CasosCriticos, MuertesHoy, Muertes)
# ---------------------------------------------------------------------------
$statUpdates = @{
    4   = @(5407564, 47262, 2831072, 2406365, 0, 996,  170127)
    5   = @(3224876, 54402, 2356640, 762773,  0, 1200, 105463)
    8   = @(572865,  3946,  437617,  123978,  0, 260,  11270)
    53  = @(45726,   462,   42180,   3379,    0, 1,    167)
    119 = @(3174,    46,    2525,    560,     0, 1,    89)
    126 = @(2597,    15,    1979,    493,     0, 0,    125)
    130 = @(2200,    11,    1558,    634,     0, 0,    8)
    137 = @(1847,    6,     949,     370,     0, 0,    528)
    138 = @(1815,    53,    577,     1158,    0, 0,    80)
    157 = @(949,     0,     860,     13,      0, 0,    76)
}

foreach ($row in $statUpdates.Keys) {
    $vals = $statUpdates[$row]
    $ws.Range("B$row:H$row").Value = $vals
}

# ---------------------------------------------------------------------------
# 3) Botsuana / Niger / Togo block (rows 150-153): Botsuana moves up to sit
#    right after Burkina Faso, Niger and Togo shift down one row, and the
#    stats for all four rows take on their refreshed values.
# ---------------------------------------------------------------------------
$ws.Range("A150:H150").Value = @("Burkina Faso", 1228, 15,  997,  177,  0, 0, 54)
$ws.Range("A151:H151").Value = @("Botsuana",     1214, 148, 120,  1091, 0, 1, 3)
$ws.Range("A152:H152").Value = @("Niger",        1161, 0,   1075, 17,   0, 0, 69)
$ws.Range("A153:H153").Value = @("Togo",         1104, 12,  791,  287,  0, 0, 26)

# ---------------------------------------------------------------------------
# 4) Montserrat / Islas Malvinas block (rows 213-214): Montserrat moves up to
#    sit right after San Bartolome, Islas Malvinas shifts down one row.
# ---------------------------------------------------------------------------
$ws.Range("A213:H213").Value = @("Montserrat",     13, 0, 12, 0, 0, 0, 1)
$ws.Range("A214:H214").Value = @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
